$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.160.46'
$ws.Range("E2").Value = '  -0.44%  '

# Row 3
$ws.Range("D3").Value = '3.330.99'
$ws.Range("E3").Value = '  -0.15%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.00%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.70%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").Value = '3.327.27'
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("E9").Value = '  -2.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.180'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.31%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.580'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.07%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.91%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '676.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.09%  '

# Row 15
$ws.Range("D15").Value = '3.862.77'
$ws.Range("E15").Value = '  -0.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.84%  '

# Row 17
$ws.Range("D17").Value = '66.294.75'
$ws.Range("E17").Value = '  -0.27%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.117'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.57%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.52%  '

# Row 20
$ws.Range("D20").Value = '3.333.38'
$ws.Range("E20").Value = '  -0.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.82%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.897'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.24%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.67%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.83%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.68%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.86%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '607.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.10%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.08%  '

# Row 35
$ws.Range("E35").Value = '  -1.46%  '

# Row 36
$ws.Range("D36").Value = '3.808.85'
$ws.Range("E36").Value = '  +2.30%  '

# Row 37
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.58%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.127'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.24%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.30%  '

# Row 41
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0700'
$ws.Range("E41").Value = '  -4.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.18%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '32.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.14%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.74%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0414'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.38%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.12%  '

# Row 48
$ws.Range("E48").Value = '  -1.85%  '

# Row 49
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.33%  '

# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.75%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.63%  '
